$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 555
$ws.Range("I33").Value = 600.1818
$ws.Range("J33").Value = 455.6
$ws.Range("K33").Value = 600.1818
$ws.Range("L33").Value = 455.6
$ws.Range("M33").Value = -371.1818
$ws.Range("N33").Value = -913.6

$ws.Range("H40").Value = 2171.5715
$ws.Range("I40").Value = 1960.2
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 1960.2
$ws.Range("L40").Value = 2700
$ws.Range("M40").Value = -1785.2
$ws.Range("N40").Value = -3050

$ws.Range("H97").Value = 404.75
$ws.Range("J97").Value = 373
$ws.Range("L97").Value = 1119
$ws.Range("N97").Value = -2111

$ws.Range("H98").Value = 4097.6665
$ws.Range("J98").Value = 1474
$ws.Range("L98").Value = 1474
$ws.Range("N98").Value = -4470

$ws.Range("H122").Value = 4097.6665
$ws.Range("J122").Value = 1474
$ws.Range("L122").Value = 4422
$ws.Range("N122").Value = -9322

$ws.Range("H129").Value = 777
$ws.Range("J129").Value = 864.875
$ws.Range("L129").Value = 2594.625
$ws.Range("N129").Value = -12594.625

$ws.Range("H137").Value = 1159.6086
$ws.Range("I137").Value = 919.1905
$ws.Range("K137").Value = 2757.5715
$ws.Range("M137").Value = -207.5715

$ws.Range("H138").Value = 1598.56
$ws.Range("I138").Value = 907.1539
$ws.Range("J138").Value = 1841.4865
$ws.Range("K138").Value = 2721.4617
$ws.Range("L138").Value = 5524.4595
$ws.Range("M138").Value = 2418.5383
$ws.Range("N138").Value = -15804.4595

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2695.4849
$ws.Range("I32").Value = 2339.9673
$ws.Range("J32").Value = 7368
$ws.Range("K32").Value = 2339.9673
$ws.Range("L32").Value = 7368
$ws.Range("M32").Value = -2052.9673
$ws.Range("N32").Value = -7942

$ws.Range("H45").Value = 1308.579
$ws.Range("I45").Value = 1155.25
$ws.Range("K45").Value = 1155.25
$ws.Range("M45").Value = -778.25

$ws.Range("H61").Value = 71429840
$ws.Range("I61").Value = 90910050
$ws.Range("K61").Value = 90910050
$ws.Range("M61").Value = -90909838

$ws.Range("H74").Value = 2408.9666
$ws.Range("I74").Value = 2018.9584
$ws.Range("K74").Value = 2018.9584
$ws.Range("M74").Value = -1144.9584

$ws.Range("H77").Value = 2408.9666
$ws.Range("I77").Value = 2018.9584
$ws.Range("K77").Value = 10094.792
$ws.Range("M77").Value = -5726.791999999999

$ws.Range("H132").Value = 2879.162
$ws.Range("I132").Value = 2408.3
$ws.Range("K132").Value = 7224.900000000001
$ws.Range("M132").Value = -4694.900000000001

$ws.Range("H136").Value = 71429840
$ws.Range("I136").Value = 90910050
$ws.Range("K136").Value = 272730150
$ws.Range("M136").Value = -272727600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 863.0909
$ws.Range("I80").Value = 331.66666
$ws.Range("J80").Value = 1062.375
$ws.Range("K80").Value = 331.66666
$ws.Range("L80").Value = 1062.375
$ws.Range("M80").Value = 666.33334
$ws.Range("N80").Value = -3058.375

$ws.Range("H83").Value = 863.0909
$ws.Range("I83").Value = 331.66666
$ws.Range("J83").Value = 1062.375
$ws.Range("K83").Value = 1658.3333
$ws.Range("L83").Value = 5311.875
$ws.Range("M83").Value = 3333.6667
$ws.Range("N83").Value = -15295.875

$ws.Range("H134").Value = 1243.4546
$ws.Range("I134").Value = 1117.85
$ws.Range("K134").Value = 3353.55
$ws.Range("M134").Value = -818.5499999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 787
$ws.Range("I107").Value = 430.2353
$ws.Range("K107").Value = 430.2353
$ws.Range("M107").Value = 1489.7647

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4127.8335
$ws.Range("J39").Value = 4193.4
$ws.Range("L39").Value = 12580.2
$ws.Range("N39").Value = -13168.2

$ws.Range("H98").Value = 1399.5555
$ws.Range("I98").Value = 1922.8334
$ws.Range("J98").Value = 353
$ws.Range("K98").Value = 5768.5002
$ws.Range("L98").Value = 1059
$ws.Range("M98").Value = -4270.5002
$ws.Range("N98").Value = -4055

$ws.Range("H107").Value = 3816.2068
$ws.Range("I107").Value = 496
$ws.Range("J107").Value = 4507.9165
$ws.Range("K107").Value = 1488
$ws.Range("L107").Value = 13523.7495
$ws.Range("M107").Value = 432
$ws.Range("N107").Value = -17363.7495

$ws.Range("H131").Value = 22225304
$ws.Range("J131").Value = 3930.5588
$ws.Range("L131").Value = 11791.6764
$ws.Range("N131").Value = -21871.6764

$ws.Range("H134").Value = 2799.5217
$ws.Range("I134").Value = 2683.7693
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 8051.3079
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -2981.3079
$ws.Range("N134").Value = -18990

$ws.Range("H136").Value = 2830.6924
$ws.Range("I136").Value = 2485.7144
$ws.Range("K136").Value = 7457.1432
$ws.Range("M136").Value = -2357.1432

$ws.Range("H138").Value = 2517.0244
$ws.Range("I138").Value = 2497.7856
$ws.Range("J138").Value = 2527
$ws.Range("K138").Value = 7493.3568
$ws.Range("L138").Value = 7581
$ws.Range("M138").Value = -2353.3568
$ws.Range("N138").Value = -17861

$ws.Range("H140").Value = 22415
$ws.Range("J140").Value = 3499.0386
$ws.Range("L140").Value = 10497.1158
$ws.Range("N140").Value = -20857.1158

$ws.Range("H141").Value = 100002830
$ws.Range("I141").Value = 111112820
$ws.Range("J141").Value = 12933
$ws.Range("K141").Value = 333338460
$ws.Range("L141").Value = 38799
$ws.Range("M141").Value = -333333280
$ws.Range("N141").Value = -49159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10100
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -8796

$ws.Range("H83").Value = 10100
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 100000
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -95008
$ws.Range("N83").Value = -43984

$ws.Range("H132").Value = 2715.4856
$ws.Range("I132").Value = 3197.0625
$ws.Range("J132").Value = 2309.9473
$ws.Range("K132").Value = 9591.1875
$ws.Range("L132").Value = 6929.841899999999
$ws.Range("M132").Value = -7061.1875
$ws.Range("N132").Value = -11989.8419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1634
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 1634
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H46").Value = 4471.4287
$ws.Range("J46").Value = 4753.846
$ws.Range("L46").Value = 4753.846
$ws.Range("N46").Value = -5129.846

$ws.Range("H61").Value = 1006.84
$ws.Range("I61").Value = 889.3889
$ws.Range("J61").Value = 1308.8572
$ws.Range("K61").Value = 889.3889
$ws.Range("L61").Value = 1308.8572
$ws.Range("M61").Value = -687.3889
$ws.Range("N61").Value = -1712.8572

$ws.Range("H68").Value = 1264.8334
$ws.Range("I68").Value = 1271.6364
$ws.Range("K68").Value = 1271.6364
$ws.Range("M68").Value = -522.6364000000001

$ws.Range("H71").Value = 1264.8334
$ws.Range("I71").Value = 1271.6364
$ws.Range("K71").Value = 6358.182000000001
$ws.Range("M71").Value = -2614.182000000001

$ws.Range("H82").Value = 2132.3333
$ws.Range("I82").Value = 2098.2
$ws.Range("J82").Value = 2200.6
$ws.Range("K82").Value = 2098.2
$ws.Range("L82").Value = 2200.6
$ws.Range("M82").Value = -1737.2
$ws.Range("N82").Value = -2922.6

$ws.Range("H85").Value = 2132.3333
$ws.Range("I85").Value = 2098.2
$ws.Range("J85").Value = 2200.6
$ws.Range("K85").Value = 2098.2
$ws.Range("L85").Value = 2200.6
$ws.Range("M85").Value = -850.1999999999998
$ws.Range("N85").Value = -4696.6

$ws.Range("H113").Value = 1006.84
$ws.Range("I113").Value = 889.3889
$ws.Range("J113").Value = 1308.8572
$ws.Range("K113").Value = 889.3889
$ws.Range("L113").Value = 1308.8572
$ws.Range("M113").Value = 1280.6111
$ws.Range("N113").Value = -5648.8572

$ws.Range("H139").Value = 52166.5
$ws.Range("J139").Value = 52166.5
$ws.Range("L139").Value = 52166.5
$ws.Range("N139").Value = -62446.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 500000000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 500000000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H97").Value = 20572
$ws.Range("J97").Value = 20572
$ws.Range("L97").Value = 20572
$ws.Range("N97").Value = -22554

$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws.Range("H132").Value = 1735.6052
$ws.Range("I132").Value = 1655.875
$ws.Range("J132").Value = 2160.8333
$ws.Range("K132").Value = 4967.625
$ws.Range("L132").Value = 6482.499899999999
$ws.Range("M132").Value = -2437.625
$ws.Range("N132").Value = -11542.4999

$ws.Range("H136").Value = 1613.5
$ws.Range("I136").Value = 1491.1
$ws.Range("J136").Value = 1817.5
$ws.Range("K136").Value = 4473.299999999999
$ws.Range("L136").Value = 5452.5
$ws.Range("M136").Value = -1923.299999999999
$ws.Range("N136").Value = -10552.5

Write-Host "Applied all edits"